# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# Row 2: 6412 -> 6434
# Row 3: 29   -> 30
# Row 5: 1011 -> 1014
# Row 6: 113  -> 114

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 6434
    $ws.Range("F3").Value = 30
    $ws.Range("F5").Value = 1014
    $ws.Range("F6").Value = 114
}
